$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-09 Monday" "2024-12-10 Tuesday"

Replace-Text "704÷8=88, 0" "547÷2=273, 1"
Replace-Text "523÷7=74, 5" "814÷6=135, 4"
Replace-Text "935÷4=233, 3" "982÷9=109, 1"
Replace-Text "391÷7=55, 6" "929÷6=154, 5"
Replace-Text "439÷3=146, 1" "306÷3=102, 0"

Replace-Text "943÷3=314, 1" "202÷8=25, 2"
Replace-Text "743÷8=92, 7" "260÷8=32, 4"
Replace-Text "533÷9=59, 2" "125÷7=17, 6"
Replace-Text "200÷8=25, 0" "981÷6=163, 3"
Replace-Text "607÷5=121, 2" "685÷2=342, 1"

Replace-Text "234÷6=39, 0" "478÷5=95, 3"
Replace-Text "134÷9=14, 8" "105÷5=21, 0"
Replace-Text "476÷3=158, 2" "336÷4=84, 0"
Replace-Text "178÷3=59, 1" "557÷8=69, 5"
Replace-Text "579÷7=82, 5" "914÷7=130, 4"

Replace-Text "723÷9=80, 3" "776÷8=97, 0"
Replace-Text "925÷4=231, 1" "773÷5=154, 3"
Replace-Text "280÷7=40, 0" "953÷9=105, 8"
Replace-Text "165÷3=55, 0" "212÷3=70, 2"
Replace-Text "117÷3=39, 0" "247÷6=41, 1"

Replace-Text "196÷8=24, 4" "107÷7=15, 2"
Replace-Text "886÷3=295, 1" "868÷8=108, 4"
Replace-Text "235÷5=47, 0" "447÷5=89, 2"
Replace-Text "911÷5=182, 1" "153÷2=76, 1"
Replace-Text "485÷4=121, 1" "139÷6=23, 1"

Write-Output "done"
